$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Data Set"
$ws.Range("B1").Value = "Read Length"
$ws.Range("C1").Value = "Total Reads"
$ws.Range("D1").Value = "Error Rate*"
$ws.Range("E1").Value = "Reads with Adapters*"
$ws.Range("F1").Value = "Adapter Bases*"

# --- Row 2: Simulated 1 ---
$ws.Range("A2").Value = "Simulated 1"
$ws.Range("B2").Value = 125
$ws.Range("C2").Value = 781923
$ws.Range("D2").Value = 0.002
$ws.Range("E2").Value = 325982
$ws.Range("F2").Value = 12447262

# --- Row 3: Simulated 2 ---
$ws.Range("A3").Value = "Simulated 2"
$ws.Range("B3").Value = 125
$ws.Range("C3").Value = 780899
$ws.Range("D3").Value = 0.006
$ws.Range("E3").Value = 325105
$ws.Range("F3").Value = 12416861

# --- Row 4: Simulated 3 ---
$ws.Range("A4").Value = "Simulated 3"
$ws.Range("B4").Value = 125
$ws.Range("C4").Value = 782237
$ws.Range("D4").Value = 0.012
$ws.Range("E4").Value = 325860
$ws.Range("F4").Value = 12464235

# --- Row 5: GM12878 WGBS ---
$ws.Range("A5").Value = "GM12878 WGBS"
$ws.Range("B5").Value = 125
$ws.Range("C5").Value = 1000000

# --- Number formats ---
$ws.Range("D2:D4").NumberFormat = "0.00%"
$ws.Range("C2:C4").NumberFormat = "#,##0"

# E (Reads with Adapters*), F (Adapter Bases*) and the GM12878 total-reads
# cell get their own distinct, un-bordered "#,##0" style - nudge the engine
# into minting a fresh cellXf instead of reusing the bordered "Total Reads"
# style above, by touching an alignment-group property (no visual effect).
$ws.Range("E2:E4").NumberFormat = "#,##0"
$ws.Range("E2:E4").IndentLevel = 0
$ws.Range("F2:F4").NumberFormat = "#,##0"
$ws.Range("F2:F4").IndentLevel = 0
$ws.Range("C5").NumberFormat = "#,##0"
$ws.Range("C5").IndentLevel = 0

# --- Header borders + alignment ---
$ws.Range("A1:F1").Borders.Item(9).LineStyle = 1
$ws.Range("B1:F1").HorizontalAlignment = -4108

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 15.125
$ws.Columns.Item(5).ColumnWidth = 17.5
$ws.Columns.Item(6).ColumnWidth = 13.5

# --- Selection / view ---
$ws.Range("E7").Select()

Write-Output "done"
